# Insert a new price-report row for "Macroferia Regional de Talca - Haba"
# (weekly refresh of the consolidated Fruta/Hortaliza data set).
#
# The new observation is inserted as row 26: we shift every row currently
# at/after row 26 down by one (Excel's native Insert behaviour), seed the
# freshly-inserted row 26 with a copy of the row that lands on top of it
# (row 27, after the shift - i.e. what used to be row 26), and then
# overwrite the three fields that make this a distinct new record: the
# report date, the traded volume, and the origin region (price columns
# and the derived $/Kg stay the same as the neighbouring record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 26..146 down to 27..147, leaving a blank row 26.
$ws.Rows("26").Insert()

# Seed the new row 26 with the values now sitting in row 27 (18 columns, A:R).
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(26, $c).Value2 = $ws.Cells.Item(27, $c).Value2
}

# Overwrite the fields that differentiate the new record.
$ws.Cells.Item(26, 4).Value2 = 45243                 # D26 - Fecha
$ws.Cells.Item(26, 10).Value2 = 400                  # J26 - Volumen
$ws.Cells.Item(26, 15).Value2 = "Región del Maule"   # O26 - Origen
